# upmachine.xlsx edit: insert a new "所在地" (location) column between
# "时间" (time) and "备注" (remarks), i.e. the new column becomes H and the
# former column H ("备注") shifts right to become column I.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Preserve the existing last-column header ("备注") and move it one column
# to the right (H -> I), then drop the new header into the vacated H cell.
$lastHeader = $ws.Cells.Item(1, 8).Value2
$ws.Cells.Item(1, 9).Value2 = $lastHeader
$ws.Cells.Item(1, 8).Value2 = "所在地"

# Move the active selection to H2, matching the saved view state.
$ws.Range("H2").Select() | Out-Null

# Nudge the shared tab-split ratio, mirroring the tiny view-state drift
# recorded in the saved workbook (best effort; cosmetic only).
$wb.Windows.Item(1).TabRatio = 0.991
